# edit.ps1 - applies the "Add Period to Figure 1, Radar Chart Section,
# References to Black" commit to the poster deck.
#
# Strategy: locate each run by its *current* full text (via
# TextRange.Characters(start,len)) and overwrite that exact span with the
# complete new text for the run. Replacing a whole run's span (rather than
# just the changed word) keeps the underlying OOXML as a single <a:r> run,
# matching how the real edit was authored.

function Replace-RunText($TextRange, $OldText, $NewText) {
    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Could not find run text: [$OldText]"
    }
    $startPos = $idx + 1
    $len = $OldText.Length
    $sub = $TextRange.Characters($startPos, $len)
    $sub.Text = $NewText
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "Inhaltsplatzhalter 9" (Focus / Symbol MAP copy) ---
$shape1 = $s.Shapes.Item(3)
$tr1 = $shape1.TextFrame.TextRange

Replace-RunText $tr1 `
    "Germany, the Netherlands and the United Kingdom, along with a line chart of World Health Organization (WHO) data on COVID-19 infection rates with an orienting vertical line. The latter contains a radar chart allowing the visitor to discover differences between all 15 search terms across the countries, as well as a small multiple time series graph displaying the changes over the study time period." `
    "Germany, the Netherlands and the United Kingdom, along with a line chart of World Health Organization (WHO) data on COVID-19 infection rates with an orienting vertical line. The latter contains a radar chart allowing the visitor to discover differences between all 15 search terms across the countries, as well as a small multiple time series plot displaying the changes over the study time period."

# --- Shape "Inhaltsplatzhalter 11" (Radar chart / Data collection copy) ---
$shape2 = $s.Shapes.Item(4)
$tr2 = $shape2.TextFrame.TextRange

Replace-RunText $tr2 `
    "Clicking on the navigation button COUNTRY a radar chart pops up comparing the 15 search terms across the countries selected in the dropdown menu (Fig. 3). A timeline animation can be started, or the time can be selected individually. N" `
    "Clicking on the navigation button COUNTRY a radar chart pops up comparing the 15 search terms across the countries selected in the dropdown menu (Fig. 3). A timeline animation can be started, or the date can be selected individually. N"

Replace-RunText $tr2 `
    " and positive values indicate lower or higher searches compared to the year before. Additional support is given by line plots presenting the evolution of each search term" `
    " and positive values indicate lower or higher searches compared to the year before. Additional context is given by line plots presenting the evolution of each search term"

# Same shape also holds the "Data collection" paragraph with the date range.
Replace-RunText $tr2 `
    "We used the Google Trends data from 1/2019 to 11/2020 and calculated the differences between the first and the second year " `
    "We used the Google Trends data from January 2019 to November 2020 and calculated the differences between the first and the second year "

# --- Shape "Rechteck 22" (Fig. 1 caption) : add trailing period + shrink box ---
$shape5 = $s.Shapes.Item(13)
$tr5 = $shape5.TextFrame.TextRange

Replace-RunText $tr5 `
    "Fig. 1 Search trends visualized by icon buttons" `
    "Fig. 1 Search trends visualized by icon buttons."

$shape5.Height = 276999 / 12700
